$d = $word.ActiveDocument

# The first paragraph in the document holds the title text that needs
# to be restyled and replaced.
$p = $d.Paragraphs.Item(1)

# Change paragraph style from "Title" to "Heading1" and center it.
$p.set_Style("Heading1")
$p.Format.Alignment = 1  # wdAlignParagraphCenter

# Replace the paragraph's text (but not its paragraph mark) with the new
# table caption, then format just that text run with black color and
# 14pt (sz=28 half-points) size.
$range = $p.Range
[void]$range.MoveEnd(1, -1)  # wdCharacter: exclude the trailing paragraph mark
$range.Text = "Таблица № 3.7. Суммарные выбросы загрязняющих веществ в атмосферу, их очистка и утилизация (в целом по предприятию), т/год."
$range.Font.Color = 0            # wdColorBlack (RGB 0,0,0 -> 000000)
$range.Font.Size = 14            # 14pt => w:sz 28 (half-points)
